# Apply updated odds values to Sheet1 of the workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("G2").Value = 1.86
$ws.Range("O2").Value = 1.39
$ws.Range("Q2").Value = 2.14
$ws.Range("S2").Value = 4
$ws.Range("T2").Value = 1.98
$ws.Range("U2").Value = 1.92
$ws.Range("W2").Value = 2.16

# Row 3
$ws.Range("L3").Value = 1.46
$ws.Range("O3").Value = 1.41
$ws.Range("P3").Value = 1.75
$ws.Range("Q3").Value = 2.3
$ws.Range("R3").Value = 1.28
$ws.Range("S3").Value = 4.2
$ws.Range("Y3").Value = 7.2
$ws.Range("AD3").Value = 10.5

# Row 4
$ws.Range("AB4").Value = 14.5
$ws.Range("AC4").Value = 8.4
$ws.Range("AF4").Value = 44
$ws.Range("AH4").Value = 28
$ws.Range("AJ4").Value = 200
$ws.Range("AL4").Value = 130
$ws.Range("AM4").Value = 250
$ws.Range("AO4").Value = 18

# Row 5
$ws.Range("N5").Value = 3.6
$ws.Range("P5").Value = 1.89
$ws.Range("T5").Value = 1.95
$ws.Range("Z5").Value = 10.5
$ws.Range("AB5").Value = 17.5
$ws.Range("AC5").Value = 9
$ws.Range("AD5").Value = 10
$ws.Range("AI5").Value = 110

# Row 7
$ws.Range("F7").Value = 3.2
$ws.Range("G7").Value = 3.3
$ws.Range("H7").Value = 2.4
$ws.Range("I7").Value = 2.44
$ws.Range("P7").Value = 1.96
$ws.Range("R7").Value = 1.37
$ws.Range("T7").Value = 1.77
$ws.Range("AB7").Value = 13
$ws.Range("AG7").Value = 14.5
$ws.Range("AI7").Value = 40
$ws.Range("AL7").Value = 48
$ws.Range("AN7").Value = 38
$ws.Range("AO7").Value = 21

# Row 8
$ws.Range("F8").Value = 1.66
$ws.Range("G8").Value = 1.68
$ws.Range("I8").Value = 6.4
$ws.Range("J8").Value = 4.2
$ws.Range("K8").Value = 4.6
$ws.Range("P8").Value = 2.16
$ws.Range("Q8").Value = 1.74
$ws.Range("X8").Value = 19.5
$ws.Range("AC8").Value = 10.5
$ws.Range("AD8").Value = 24
$ws.Range("AE8").Value = 80
$ws.Range("AG8").Value = 11.5
$ws.Range("AI8").Value = 75
$ws.Range("AJ8").Value = 16.5
$ws.Range("AK8").Value = 17.5
$ws.Range("AN8").Value = 8.4
